$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.784.28'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '3.493.45'
$ws.Range("E3").Value = '  -1.72%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.51'
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '198.63'
$ws.Range("E6").Value = '  +7.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("E7").Value = '  +1.61%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.210'
$ws.Range("E9").Value = '  -1.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.654'
$ws.Range("E10").Value = '  +1.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.25'
$ws.Range("E11").Value = '  +1.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000304'
$ws.Range("E12").Value = '  -1.23%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.56'
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("D14").Value = '4.055.24'
$ws.Range("E14").Value = '  -1.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '595.91'
$ws.Range("E15").Value = '  +3.24%  '
$ws.Range("D16").Value = '69.866.76'
$ws.Range("E16").Value = '  -0.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.95'
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.63'
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("D19").Value = '3.479.70'
$ws.Range("E19").Value = '  -2.50%  '
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.987'
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.90'
$ws.Range("E22").Value = '  +2.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '103.46'
$ws.Range("E23").Value = '  +9.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.64'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.04'
$ws.Range("E25").Value = '  +3.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.11'
$ws.Range("E26").Value = '  +6.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.94'
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.83'
$ws.Range("E28").Value = '  +5.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.53'
$ws.Range("E29").Value = '  +4.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.55'
$ws.Range("E30").Value = '  +23.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.24'
$ws.Range("E31").Value = '  +3.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.77'
$ws.Range("E32").Value = '  +4.76%  '
$ws.Range("E33").Value = '  +1.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.69'
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").Value = '3.706.19'
$ws.Range("E35").Value = '  +4.25%  '
$ws.Range("B36").Value = 'PEPE'
$ws.Range("C36").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D36").Value = '0.0₃0808'
$ws.Range("E36").Value = '  +3.64%  '
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '521.72'
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.20%  '
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.391'
$ws.Range("E39").Value = '  -3.20%  '
$ws.Range("B40").Value = 'InjectiveProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.80'
$ws.Range("E40").Value = '  -1.27%  '
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.98'
$ws.Range("E41").Value = '  -6.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.55'
$ws.Range("E42").Value = '  +0.67%  '
$ws.Range("E43").Value = '  +0.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0458'
$ws.Range("E44").Value = '  +0.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.85'
$ws.Range("E45").Value = '  -1.99%  '
$ws.Range("E46").Value = '  -0.30%  '
$ws.Range("E47").Value = '  -3.72%  '
$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.01'
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.75'
$ws.Range("E49").Value = '  -4.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.12'
$ws.Range("E50").Value = '  -2.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000242'
$ws.Range("E51").Value = '  -0.66%  '

Write-Host "Applied cryptos update"